# Generate Report for Archive
#
# - Flip the localization status from "Ready for handoff" to
#   "In Translation" everywhere it is shown:
#     Overview!E2:F3, zh-cn!C2:C3, de-de!C2:C3
# - Narrow the now-shorter status columns (Overview E:F, zh-cn C, de-de C)
#   to match the new, narrower column width.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn     = $wb.Worksheets.Item("zh-cn")
$dede     = $wb.Worksheets.Item("de-de")

# --- Status text: "Ready for handoff" -> "In Translation" ---
$overview.Range("E2:F3").Value = "In Translation"
$zhcn.Range("C2:C3").Value = "In Translation"
$dede.Range("C2:C3").Value = "In Translation"

# --- Narrow the status columns to match the new width ---
$overview.Range("E1:F1").ColumnWidth = 12.5
$zhcn.Range("C1").ColumnWidth = 12.5
$dede.Range("C1").ColumnWidth = 12.5
